# Round 23 (and final) results + tips tidy-up for the AFL 2024 fixture sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fixture")

# Results (column G) and tip columns (H:AE) for match rows 182-199
# (Round 22 tail + Round 23), which previously only had the
# Match Number / Round / Date / Location / Home / Away columns filled in.
$data = @{}
$data[182] = @("89 - 86", 1, 1, 0, 1, 1, 1, 1, 1, 1, 1, 1, 0, 1, 0, 1, 1, 1, 1, -1, 1, 1, 1, 1, 1)
$data[183] = @("64 - 82", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, -1, 1, 1, 1, 1, 1)
$data[184] = @("97 - 102", 1, 1, 0, 1, 1, 1, 1, 1, 1, 1, 1, 0, 0, 1, 1, 1, 1, 1, -1, 1, 1, 1, 1, 1)
$data[185] = @("62 - 73", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 0, 0, 0, 1, 1, 1, -1, 1, 1, 1, 1, 1)
$data[186] = @("86 - 87", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 1, -1, 1, 1, 1, 1, 1)
$data[187] = @("51 - 53", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, -1, 0, 0, 0, 0, 0)
$data[188] = @("38 - 112", 1, 0, 1, 1, 1, 1, 1, 0, 0, 1, 1, 0, 1, 1, 1, 0, 1, 1, -1, 1, 1, 0, 0, 1)
$data[189] = @("51 - 99", 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, -1, 0, 0, 0, 0, 0)
$data[190] = @("111 - 72", 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, -1, 0, 0, 0, 0, 0)
$data[191] = @("59 - 98", 1, 1, 1, 0, 1, 0, 1, 0, 0, 1, 0, 0, 0, 0, 1, 0, 0, 0, -1, 0, 0, 1, 0, 0)
$data[192] = @("63 - 117", 0, 1, 0, 1, 1, 1, 1, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 1, -1, 1, 1, 1, 1, 1)
$data[193] = @("101 - 92", 0, 0, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 1, -1, 1, 1, 0, 1, 1)
$data[194] = @("79 - 78", 1, 0, 1, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, -1, 0, 0, 1, 0, 0)
$data[195] = @("107 - 89", 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, -1, 0, 0, 0, 0, 0)
$data[196] = @("80 - 58", 1, 1, 0, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 0, 1, 1, 1, 0, -1, 1, 1, 0, 1, 1)
$data[197] = @("138 - 42", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, -1, 1, 1, 1, 1, 1)
$data[198] = @("131 - 68", 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, -1, 1, 1, 1, 1, 1)
$data[199] = @("34 - 99", 0, 0, 1, 1, 1, 1, 0, 1, 0, 0, 1, 1, 1, 0, 1, 0, 1, 0, -1, 1, 0, 1, 1, 1)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 7 + $i).Value = $vals[$i]
    }
}

# Widen the "margin" column (T) slightly to fit the new larger numbers.
$ws.Columns.Item(20).ColumnWidth = 7.5703125

# Move the frozen-pane view / active selection down to where the
# newly-entered round 23 results now sit.
$ws.Application.ActiveWindow.ScrollRow = 171
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("AE200").Select()
